# Renaming of ecology_format tables
#
# "ecological_params" sheet header (B1:D1) is renamed:
#   PP_impellar_virgin_(recipe_endpoint_h)         -> PP_virgin
#   PP_impellar_recycled_(recipe_endpoint_h)       -> PP_recycled
#   PP_impellar_recycled_vision_(recipe_endpoint_h)-> PP_recycled_industrial
# and its header formatting is switched to match the "Scaling" sheet's
# plain (non-bold, borderless) header style.
#
# The active/selected sheet also moves from "Scaling" to "ecological_params",
# with a new cell selection on each sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ecological_params")
$ws2 = $wb.Worksheets.Item("Scaling")

# Rename the three header labels.
$ws1.Range("B1").Value = "PP_virgin"
$ws1.Range("C1").Value = "PP_recycled"
$ws1.Range("D1").Value = "PP_recycled_industrial"

# Match the header formatting used on the "Scaling" sheet (plain font,
# no border) by copying its B1:D1 formats onto the renamed headers.
$ws2.Range("B1:D1").Copy()
$ws1.Range("B1:D1").PasteSpecial(-4122)

# "ecological_params" becomes the active sheet/selection; "Scaling" keeps
# its own (already-correct) selection but loses the active-tab flag.
$ws1.Activate()
$ws1.Range("D16").Select()
